# Updated cryptos list on Thu Jul 11 21:11:53 UTC 2024 with GitHub Actions
#
# Refresh the live price / 1h-volume-change snapshot for each coin row.
# Numeric-looking Price values are written with a leading apostrophe so
# Excel keeps storing them as text (matching the sheet's existing
# inlineStr cells) instead of silently re-typing them as numbers; the
# Style reset at the end clears the resulting "quote prefix" number
# format again so no visible formatting change is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.553.94'
$ws.Range("E2").Value = '  +0.22%  '

$ws.Range("D3").Value = '3.115.57'
$ws.Range("E3").Value = '  +0.74%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = "'526.47"
$ws.Range("E5").Value = '  +0.82%  '

$ws.Range("D6").Value = "'137.32"
$ws.Range("E6").Value = '  -2.47%  '

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").Value = '3.113.93'
$ws.Range("E8").Value = '  +0.71%  '

$ws.Range("E9").Value = '  +2.44%  '

$ws.Range("D10").Value = "'7.26"
$ws.Range("E10").Value = '  +0.94%  '

$ws.Range("E11").Value = '  -0.03%  '

$ws.Range("D12").Value = "'0.395"
$ws.Range("E12").Value = '  +2.95%  '

$ws.Range("D13").Value = '3.652.18'
$ws.Range("E13").Value = '  +0.71%  '

$ws.Range("E14").Value = '  +2.99%  '

$ws.Range("D15").Value = "'25.26"
$ws.Range("E15").Value = '  -2.65%  '

$ws.Range("D16").Value = "'0.0000163"
$ws.Range("E16").Value = '  +0.26%  '

$ws.Range("D17").Value = '57.645.64'
$ws.Range("E17").Value = '  +0.27%  '

$ws.Range("D18").Value = '3.113.89'
$ws.Range("E18").Value = '  +0.49%  '

$ws.Range("D19").Value = "'5.94"
$ws.Range("E19").Value = '  -2.48%  '

$ws.Range("D20").Value = "'12.42"
$ws.Range("E20").Value = '  -2.65%  '

$ws.Range("E21").Value = '  -1.58%  '

$ws.Range("D22").Value = "'348.42"
$ws.Range("E22").Value = '  +2.82%  '

$ws.Range("D23").Value = "'5.80"
$ws.Range("E23").Value = '  -0.32%  '

$ws.Range("E24").Value = '  -0.05%  '

$ws.Range("D25").Value = "'68.04"
$ws.Range("E25").Value = '  +2.14%  '

$ws.Range("D26").Value = "'0.502"
$ws.Range("E26").Value = '  -1.83%  '

$ws.Range("D27").Value = "'0.167"
$ws.Range("E27").Value = '  -0.35%  '

$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = '  -0.60%  '

$ws.Range("D29").Value = '0.0₃0908'
$ws.Range("E29").Value = '  +0.15%  '

$ws.Range("D30").Value = "'7.45"
$ws.Range("E30").Value = '  +3.96%  '

$ws.Range("E31").Value = '  +0.06%  '

$ws.Range("E32").Value = '  +0.77%  '

$ws.Range("D33").Value = "'6.08"
$ws.Range("E33").Value = '  -6.16%  '

$ws.Range("D34").Value = "'20.98"
$ws.Range("E34").Value = '  +0.37%  '

$ws.Range("D35").Value = "'4.97"
$ws.Range("E35").Value = '  +7.66%  '

$ws.Range("E36").Value = '  -2.07%  '

$ws.Range("D37").Value = "'158.69"
$ws.Range("E37").Value = '  +1.40%  '

$ws.Range("E38").Value = '  +0.15%  '

$ws.Range("D39").Value = "'26.06"
$ws.Range("E39").Value = '  -4.10%  '

$ws.Range("D40").Value = "'1.24"
$ws.Range("E40").Value = '  -2.88%  '

$ws.Range("E41").Value = '  +7.21%  '

$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = "'1.62"
$ws.Range("E42").Value = '  +6.62%  '

$ws.Range("B43").Value = 'Hedera'
$ws.Range("C43").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D43").Value = "'0.0664"
$ws.Range("E43").Value = '  +1.10%  '

$ws.Range("D44").Value = "'0.702"
$ws.Range("E44").Value = '  +2.61%  '

$ws.Range("D45").Value = '3.152.64'
$ws.Range("E45").Value = '  +0.50%  '

$ws.Range("D46").Value = '2.356.74'
$ws.Range("E46").Value = '  +1.89%  '

$ws.Range("D47").Value = "'36.52"
$ws.Range("E47").Value = '  -0.29%  '

$ws.Range("D48").Value = "'0.999"
$ws.Range("E48").Value = '  -0.01%  '

$ws.Range("D49").Value = "'0.0268"
$ws.Range("E49").Value = '  +3.35%  '

$ws.Range("D50").Value = "'0.962"
$ws.Range("E50").Value = '  -1.03%  '

$ws.Range("E51").Value = '  +0.52%  '

# Clear the text-number-format residue left by the apostrophe-prefixed
# Price values above, restoring the cells to the sheet's default style.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
